$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.948.26'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = '3.062.33'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = "'539.46"
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('D6').Value = "'136.99"
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.20%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').Value = '3.053.54'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('E9').Value = '  +1.54%  '
$ws.Range('E10').Value = '  +1.54%  '
$ws.Range('E11').Value = '  +2.12%  '
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('E13').Value = '  +3.85%  '
$ws.Range('D14').Value = "'34.36"
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = '3.562.82'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').Value = '62.959.76'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = "'0.112"
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.93%  '
$ws.Range('D18').Value = '3.064.21'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('D19').Value = "'6.60"
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = "'467.95"
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').Value = "'13.50"
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.40%  '
$ws.Range('D22').Value = "'0.694"
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = "'7.01"
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.46%  '
$ws.Range('D24').Value = "'78.38"
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.86%  '
$ws.Range('D25').Value = "'12.07"
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  -4.73%  '
$ws.Range('D29').Value = "'0.999"
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').Value = "'26.04"
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  +5.05%  '
$ws.Range('E32').Value = '  -2.42%  '
$ws.Range('D33').Value = "'58.73"
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.27%  '
$ws.Range('D34').Value = "'2.31"
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.65%  '
$ws.Range('D35').Value = "'5.44"
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.06%  '
$ws.Range('D36').Value = "'5.94"
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('D37').Value = "'480.66"
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.90%  '
$ws.Range('D38').Value = '3.240.42'
$ws.Range('E38').Value = '  +4.20%  '
$ws.Range('E39').Value = '  +1.53%  '
$ws.Range('D40').Value = "'0.0790"
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.35%  '
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('D42').Value = "'8.11"
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.67%  '
$ws.Range('E43').Value = '  +0.77%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = "'122.89"
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.79%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = "'25.09"
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.00%  '
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('E49').Value = '  +1.80%  '
$ws.Range('E50').Value = '  +3.40%  '
$ws.Range('E51').Value = '  +0.97%  '
